# Populate the "Inputs" sheet's first client column (D) with the McLaren
# example-client seed values, and set the preliminary revenue-split
# percentages on "Outputs_Internal". All other changed cells in the
# workbook are formula-driven and will recalculate automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Inputs sheet - column D (first client)
# ---------------------------------------------------------------------
$inputs = $wb.Worksheets.Item("Inputs")

$inputs.Range("D4").Value  = "Y"
$inputs.Range("D5").Value  = "N"
$inputs.Range("D6").Value  = "N"
$inputs.Range("D7").Value  = 50000
$inputs.Range("D8").Value  = "Y"
$inputs.Range("D9").Value  = "Y"
$inputs.Range("D10").Value = "Y"
$inputs.Range("D11").Value = "Y"
$inputs.Range("D12").Value = "Y"
$inputs.Range("D13").Value = "N"
$inputs.Range("D14").Value = "N"
$inputs.Range("D15").Value = "N"
$inputs.Range("D16").Value = "N"
$inputs.Range("D17").Value = "N"
$inputs.Range("D18").Value = "Y"
$inputs.Range("D19").Value = "N"
$inputs.Range("D20").Value = "N"
$inputs.Range("D21").Value = "N"

# ---------------------------------------------------------------------
# Outputs_Internal sheet - preliminary revenue split percentages (F/G)
# ---------------------------------------------------------------------
$internal = $wb.Worksheets.Item("Outputs_Internal")

$internal.Range("F4").Value  = 0.1
$internal.Range("G4").Value  = 0.9

$internal.Range("F7").Value  = 0.3
$internal.Range("G7").Value  = 0.7

$internal.Range("F8").Value  = 0.2
$internal.Range("G8").Value  = 0.8

$internal.Range("F9").Value  = 0.1
$internal.Range("G9").Value  = 0.9

$internal.Range("F10").Value = 0.1
$internal.Range("G10").Value = 0.9

$internal.Range("F11").Value = 0.9
$internal.Range("G11").Value = 0.1

$internal.Range("F12").Value = 0.1
$internal.Range("G12").Value = 0.9

$internal.Range("F15").Value = 0.02
$internal.Range("G15").Value = 0.98

$internal.Range("F16").Value = 0.02
$internal.Range("G16").Value = 0.98

$internal.Range("F17").Value = 0.02
$internal.Range("G17").Value = 0.98

$internal.Range("F35").Value = 1
$internal.Range("G35").Value = 0

$internal.Range("F36").Value = 1
$internal.Range("G36").Value = 0

$internal.Range("F37").Value = 0.02
$internal.Range("G37").Value = 0.98

$internal.Range("F38").Value = 0.02
$internal.Range("G38").Value = 0.98

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping to match the authored workbook view
# ---------------------------------------------------------------------
$inputs.Range("D4").Select() | Out-Null
$internal.Range("H60").Select() | Out-Null

$excel.CalculateFullRebuild()

$wb.Save()
